$d = $word.ActiveDocument

# Create new "join" paragraph styles, each based on the existing MSCJoin
# style, mirroring the existing MSC_Paragraph_A/B/C pattern that is based
# on MSCParagraph.
$wdStyleTypeParagraph = 1

$joinA = $d.Styles.Add("MSC_Join_A", $wdStyleTypeParagraph)
$joinA.BaseStyle = "MSCJoin"

$joinB = $d.Styles.Add("MSC_Join_B", $wdStyleTypeParagraph)
$joinB.BaseStyle = "MSCJoin"

$joinC = $d.Styles.Add("MSC_Join_C", $wdStyleTypeParagraph)
$joinC.BaseStyle = "MSCJoin"

# Re-point every paragraph currently using the MSCJoin style onto the new
# MSC_Join_A style (the "[...]" ellipsis join paragraphs).
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "MSC_Join") {
        $p.Style = "MSC_Join_A"
    }
}

Write-Output "done"
